# Apply the data refresh for the "Pomelo" weekly logic sheet.
# The rows get re-shuffled: row 2's data moves out, and is replaced by what
# used to be row 7 (after the underlying weekly rotation); row 3<->row4 swap,
# row 6 gets what was row 2, row 7 gets what was row 8, row 8 gets what was row 6.
# Rows 5 and 9 are unaffected. Columns A,B,C,E,F,G,H,I,J,K,L,R,T are identical
# across these rows and remain unchanged; only D (Fecha), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion) and S (Precio $/Kg) change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the "before" values for the rows that change, keyed by row number.
$before = @{
    2 = @{ D = 44400; M = 100; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos";           S = 714 }
    3 = @{ D = 44397; M = 60;  N = 11000; O = 11000; P = 11000; Q = "`$/caja 14 kilos";           S = 786 }
    4 = @{ D = 44351; M = 300; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 }
    6 = @{ D = 44162; M = 120; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 }
    7 = @{ D = 44176; M = 250; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 }
    8 = @{ D = 44208; M = 210; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 }
}

# Mapping of each target (after) row to the row whose original ("before") data it now holds.
$rowSource = @{
    2 = 7
    3 = 4
    4 = 3
    6 = 2
    7 = 8
    8 = 6
}

foreach ($targetRow in $rowSource.Keys) {
    $src = $before[$rowSource[$targetRow]]

    $ws.Cells.Item($targetRow, 4).Value  = $src.D   # Column D - Fecha
    $ws.Cells.Item($targetRow, 13).Value = $src.M   # Column M - Volumen
    $ws.Cells.Item($targetRow, 14).Value = $src.N   # Column N - Precio minimo
    $ws.Cells.Item($targetRow, 15).Value = $src.O   # Column O - Precio maximo
    $ws.Cells.Item($targetRow, 16).Value = $src.P   # Column P - Precio promedio ponderado
    $ws.Cells.Item($targetRow, 17).Value = $src.Q   # Column Q - Unidad de comercializacion
    $ws.Cells.Item($targetRow, 19).Value = $src.S   # Column S - Precio $/Kg
}
